$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.739652666666667
$ws.Range("H2").Value = 23.218958
$ws.Range("I2").Value = 0.6488398532974882
$ws.Range("J2").Value = 0.6488398532974882
$ws.Range("M2").Value = 36.81180933333333
$ws.Range("N2").Value = 110.435428
$ws.Range("O2").Value = 0.2598784967371026
$ws.Range("P2").Value = 0.2598784967371026
$ws.Range("Q2").Value = 284.9106182715582
$ws.Range("R2").Value = 2564.195564444024
$ws.Range("S2").Value = 0.1686195256980734
$ws.Range("T2").Value = 0.1686195256980734

$ws.Range("G3").Value = 7.739652666666667
$ws.Range("H3").Value = 23.218958
$ws.Range("I3").Value = 0.6488398532974882
$ws.Range("J3").Value = 0.6488398532974882
$ws.Range("O3").Value = 0.1970278712683331
$ws.Range("P3").Value = 0.197027871268333
$ws.Range("Q3").Value = 216.0060694693689
$ws.Range("R3").Value = 1944.05462522432
$ws.Range("S3").Value = 0.1278395350892616
$ws.Range("T3").Value = 0.1278395350892616

$ws.Range("G4").Value = 7.739652666666667
$ws.Range("H4").Value = 23.218958
$ws.Range("I4").Value = 0.6488398532974882
$ws.Range("J4").Value = 0.6488398532974882
$ws.Range("M4").Value = 21.95609833333333
$ws.Range("N4").Value = 65.868295
$ws.Range("O4").Value = 0.1550023737603119
$ws.Range("P4").Value = 0.1550023737603119
$ws.Range("Q4").Value = 169.9325750151789
$ws.Range("R4").Value = 1529.39317513661
$ws.Range("S4").Value = 0.1005717174514032
$ws.Range("T4").Value = 0.1005717174514032

$ws.Range("G5").Value = 7.739652666666667
$ws.Range("H5").Value = 23.218958
$ws.Range("I5").Value = 0.6488398532974882
$ws.Range("J5").Value = 0.6488398532974882
$ws.Range("M5").Value = 13.23098133333333
$ws.Range("N5").Value = 39.692944
$ws.Range("O5").Value = 0.09340609987756826
$ws.Range("P5").Value = 0.09340609987756825
$ws.Range("Q5").Value = 102.4031999591502
$ws.Range("R5").Value = 921.628799632352
$ws.Range("S5").Value = 0.06060560014165191
$ws.Range("T5").Value = 0.06060560014165191

$ws.Range("G6").Value = 7.739652666666667
$ws.Range("H6").Value = 23.218958
$ws.Range("I6").Value = 0.6488398532974882
$ws.Range("J6").Value = 0.6488398532974882
$ws.Range("M6").Value = 22.080681
$ws.Range("N6").Value = 66.242043
$ws.Range("O6").Value = 0.1558818838066577
$ws.Range("P6").Value = 0.1558818838066577
$ws.Range("Q6").Value = 170.896801583466
$ws.Range("R6").Value = 1538.071214251194
$ws.Range("S6").Value = 0.1011423786208479
$ws.Range("T6").Value = 0.1011423786208479

$ws.Range("G7").Value = 7.739652666666667
$ws.Range("H7").Value = 23.218958
$ws.Range("I7").Value = 0.6488398532974882
$ws.Range("J7").Value = 0.6488398532974882
$ws.Range("M7").Value = 19.66149466666667
$ws.Range("N7").Value = 58.984484
$ws.Range("O7").Value = 0.1388032745500265
$ws.Range("P7").Value = 0.1388032745500265
$ws.Range("Q7").Value = 152.1731396275191
$ws.Range("R7").Value = 1369.558256647672
$ws.Range("S7").Value = 0.09006109629625017
$ws.Range("T7").Value = 0.09006109629625017

$ws.Range("I8").Value = 0.3053032463428815
$ws.Range("J8").Value = 0.3053032463428815
$ws.Range("M8").Value = 36.81180933333333
$ws.Range("N8").Value = 110.435428
$ws.Range("O8").Value = 0.2598784967371026
$ws.Range("P8").Value = 0.2598784967371026
$ws.Range("Q8").Value = 134.0610263592773
$ws.Range("R8").Value = 1206.549237233496
$ws.Range("S8").Value = 0.07934174870854536
$ws.Range("T8").Value = 0.07934174870854536

$ws.Range("I9").Value = 0.3053032463428815
$ws.Range("J9").Value = 0.3053032463428815
$ws.Range("O9").Value = 0.1970278712683331
$ws.Range("P9").Value = 0.197027871268333
$ws.Range("S9").Value = 0.06015324871824943
$ws.Range("T9").Value = 0.06015324871824942

$ws.Range("I10").Value = 0.3053032463428815
$ws.Range("J10").Value = 0.3053032463428815
$ws.Range("M10").Value = 21.95609833333333
$ws.Range("N10").Value = 65.868295
$ws.Range("O10").Value = 0.1550023737603119
$ws.Range("P10").Value = 0.1550023737603119
$ws.Range("Q10").Value = 79.95958717374334
$ws.Range("R10").Value = 719.6362845636901
$ws.Range("S10").Value = 0.0473227278998759
$ws.Range("T10").Value = 0.0473227278998759

$ws.Range("I11").Value = 0.3053032463428815
$ws.Range("J11").Value = 0.3053032463428815
$ws.Range("M11").Value = 13.23098133333333
$ws.Range("N11").Value = 39.692944
$ws.Range("O11").Value = 0.09340609987756826
$ws.Range("P11").Value = 0.09340609987756825
$ws.Range("Q11").Value = 48.18450843384534
$ws.Range("R11").Value = 433.660575904608
$ws.Range("S11").Value = 0.02851718552084901
$ws.Range("T11").Value = 0.02851718552084901

$ws.Range("I12").Value = 0.3053032463428815
$ws.Range("J12").Value = 0.3053032463428815
$ws.Range("M12").Value = 22.080681
$ws.Range("N12").Value = 66.242043
$ws.Range("O12").Value = 0.1558818838066577
$ws.Range("P12").Value = 0.1558818838066577
$ws.Range("Q12").Value = 80.41329158171401
$ws.Range("R12").Value = 723.719624235426
$ws.Range("S12").Value = 0.04759124517221645
$ws.Range("T12").Value = 0.04759124517221645

$ws.Range("I13").Value = 0.3053032463428815
$ws.Range("J13").Value = 0.3053032463428815
$ws.Range("M13").Value = 19.66149466666667
$ws.Range("N13").Value = 58.984484
$ws.Range("O13").Value = 0.1388032745500265
$ws.Range("P13").Value = 0.1388032745500265
$ws.Range("Q13").Value = 71.60311330809867
$ws.Range("R13").Value = 644.4280197728881
$ws.Range("S13").Value = 0.04237709032314535
$ws.Range("T13").Value = 0.04237709032314535

$ws.Range("G14").Value = 0.5470016666666667
$ws.Range("H14").Value = 1.641005
$ws.Range("I14").Value = 0.04585690035963046
$ws.Range("J14").Value = 0.04585690035963046
$ws.Range("M14").Value = 36.81180933333333
$ws.Range("N14").Value = 110.435428
$ws.Range("O14").Value = 0.2598784967371026
$ws.Range("P14").Value = 0.2598784967371026
$ws.Range("Q14").Value = 20.13612105834889
$ws.Range("R14").Value = 181.22508952514
$ws.Range("S14").Value = 0.01191722233048387
$ws.Range("T14").Value = 0.01191722233048386

$ws.Range("G15").Value = 0.5470016666666667
$ws.Range("H15").Value = 1.641005
$ws.Range("I15").Value = 0.04585690035963046
$ws.Range("J15").Value = 0.04585690035963046
$ws.Range("O15").Value = 0.1970278712683331
$ws.Range("P15").Value = 0.197027871268333
$ws.Range("Q15").Value = 15.26627680835556
$ws.Range("R15").Value = 137.3964912752
$ws.Range("S15").Value = 0.009035087460822048
$ws.Range("T15").Value = 0.009035087460822044

$ws.Range("G16").Value = 0.5470016666666667
$ws.Range("H16").Value = 1.641005
$ws.Range("I16").Value = 0.04585690035963046
$ws.Range("J16").Value = 0.04585690035963046
$ws.Range("M16").Value = 21.95609833333333
$ws.Range("N16").Value = 65.868295
$ws.Range("O16").Value = 0.1550023737603119
$ws.Range("P16").Value = 0.1550023737603119
$ws.Range("Q16").Value = 12.01002238183056
$ws.Range("R16").Value = 108.090201436475
$ws.Range("S16").Value = 0.007107928409032824
$ws.Range("T16").Value = 0.007107928409032822

$ws.Range("G17").Value = 0.5470016666666667
$ws.Range("H17").Value = 1.641005
$ws.Range("I17").Value = 0.04585690035963046
$ws.Range("J17").Value = 0.04585690035963046
$ws.Range("M17").Value = 13.23098133333333
$ws.Range("N17").Value = 39.692944
$ws.Range("O17").Value = 0.09340609987756826
$ws.Range("P17").Value = 0.09340609987756825
$ws.Range("Q17").Value = 7.237368840968889
$ws.Range("R17").Value = 65.13631956872
$ws.Range("S17").Value = 0.004283314215067339
$ws.Range("T17").Value = 0.004283314215067338

$ws.Range("G18").Value = 0.5470016666666667
$ws.Range("H18").Value = 1.641005
$ws.Range("I18").Value = 0.04585690035963046
$ws.Range("J18").Value = 0.04585690035963046
$ws.Range("M18").Value = 22.080681
$ws.Range("N18").Value = 66.242043
$ws.Range("O18").Value = 0.1558818838066577
$ws.Range("P18").Value = 0.1558818838066577
$ws.Range("Q18").Value = 12.078169308135
$ws.Range("R18").Value = 108.703523773215
$ws.Range("S18").Value = 0.007148260013593396
$ws.Range("T18").Value = 0.007148260013593395

$ws.Range("G19").Value = 0.5470016666666667
$ws.Range("H19").Value = 1.641005
$ws.Range("I19").Value = 0.04585690035963046
$ws.Range("J19").Value = 0.04585690035963046
$ws.Range("M19").Value = 19.66149466666667
$ws.Range("N19").Value = 58.984484
$ws.Range("O19").Value = 0.1388032745500265
$ws.Range("P19").Value = 0.1388032745500265
$ws.Range("Q19").Value = 10.75487035182445
$ws.Range("R19").Value = 96.79383316642
$ws.Range("S19").Value = 0.006365087930630996
$ws.Range("T19").Value = 0.006365087930630995
